$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert two new rows at row 168 (pushes existing rows 168-292 down to 170-294)
$ws.Rows("168:169").Insert()

# New week of data (row 168 = "Primera", row 169 = "Segunda")
$ws.Range("A168").Value = 1
$ws.Range("B168").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C168").Value = "Arica y Parinacota"
$ws.Range("D168").Value = 44669
$ws.Range("E168").Value = 15
$ws.Range("F168").Value = 100114014
$ws.Range("G168").Value = "Betarraga"
$ws.Range("H168").Value = "Sin especificar"
$ws.Range("I168").Value = "Primera"
$ws.Range("J168").Value = 900
$ws.Range("K168").Value = 500
$ws.Range("L168").Value = 600
$ws.Range("M168").Value = 550
$ws.Range("N168").Value = "`$/paquete 4 unidades"
$ws.Range("O168").Value = "Región de Arica y Parinacota"
$ws.Range("P168").Value = 138
$ws.Range("Q168").Value = 4
$ws.Range("R168").Value = "Hortaliza"

$ws.Range("A169").Value = 1
$ws.Range("B169").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C169").Value = "Arica y Parinacota"
$ws.Range("D169").Value = 44669
$ws.Range("E169").Value = 15
$ws.Range("F169").Value = 100114014
$ws.Range("G169").Value = "Betarraga"
$ws.Range("H169").Value = "Sin especificar"
$ws.Range("I169").Value = "Segunda"
$ws.Range("J169").Value = 800
$ws.Range("K169").Value = 500
$ws.Range("L169").Value = 600
$ws.Range("M169").Value = 550
$ws.Range("N169").Value = "`$/paquete 5 unidades"
$ws.Range("O169").Value = "Región de Arica y Parinacota"
$ws.Range("P169").Value = 110
$ws.Range("Q169").Value = 5
$ws.Range("R169").Value = "Hortaliza"

# Ensure the date cells carry the same number format as the rest of column D
$ws.Range("D168:D169").NumberFormat = $ws.Range("D170").NumberFormat
